$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -1.79%  "
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("E11").Value = "  -2.04%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -4.19%  "
$ws.Range("E30").Value = "  -11.22%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("E32").Value = "  -4.51%  "
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("E34").Value = "  -4.71%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("E42").Value = "  -6.96%  "
$ws.Range("E43").Value = "  -9.37%  "
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("E51").Value = "  -3.12%  "
